$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

$ws.Range("C14").Value = 'Archived %1$d note'
$ws.Range("J14").Value = '%1$d note archivée'
$ws.Range("M14").Value = 'Archiviata %1$d nota'
$ws.Range("Q14").Value = 'Gearchiveerde %1$d notitie'
$ws.Range("C15").Value = 'Archived %1$d notes'
$ws.Range("J15").Value = '%1$d notes archivées'
$ws.Range("M15").Value = 'Archiviate %1$d note'
$ws.Range("Q15").Value = 'Gearchiveerde %1$d notities'
$ws.Range("C39").Value = 'Can’t add %1$d file'
$ws.Range("J39").Value = '%1$d fichier n\''a pas pu être ajouté'
$ws.Range("M39").Value = 'Impossibile aggiungere %1$d file'
$ws.Range("M40").Value = 'Impossibile aggiungere %1$d file'
$ws.Range("Q39").Value = 'Kan %1$d bestand niet toevoegen'
$ws.Range("C40").Value = 'Can’t add %1$d files'
$ws.Range("J40").Value = '%1$d fichiers n\''ont pas pu être ajoutés'
$ws.Range("Q40").Value = 'Kan %1$d bestanden niet toevoegen'
$ws.Range("E43").Value = 'Nepodařilo se přidat %1$d obrázky'
$ws.Range("S43").Value = 'Nie można dodać %1$d obrazów'
$ws.Range("S44").Value = 'Nie można dodać %1$d obrazów'
$ws.Range("S46").Value = 'Nie można dodać %1$d obrazów'
$ws.Range("Y43").Value = '%1$d slike niso bile dodane.'
$ws.Range("E44").Value = 'Nepodařilo se přidat %1$d obrázků'
$ws.Range("C45").Value = 'Can’t add %1$d image'
$ws.Range("E45").Value = 'Nepodařilo se přidat %1$d obrázek'
$ws.Range("G45").Value = 'Kann %1$d Bild nicht hinzufügen'
$ws.Range("J45").Value = 'Impossible d\''ajouter %1$d image'
$ws.Range("M45").Value = 'Impossibile aggiungere %1$d immagine'
$ws.Range("P45").Value = 'Kan ikke legge til %1$d bilde'
$ws.Range("Q45").Value = 'Kan %1$d afbeelding niet toevoegen'
$ws.Range("R45").Value = 'Kan ikkje legga til %1$d bilete'
$ws.Range("R46").Value = 'Kan ikkje legga til %1$d bilete'
$ws.Range("S45").Value = 'Nie można dodać %1$d obrazu'
$ws.Range("Y45").Value = '%1$d slika ni bila dodana.'
$ws.Range("AD45").Value = 'Không thể thêm %1$d ảnh'
$ws.Range("AD46").Value = 'Không thể thêm %1$d ảnh'
$ws.Range("C46").Value = 'Can’t add %1$d images'
$ws.Range("G46").Value = 'Kann %1$d Bilder nicht hinzufügen'
$ws.Range("J46").Value = 'Impossible d\''ajouter %1$d images'
$ws.Range("M46").Value = 'Impossibile aggiungere %1$d immagini'
$ws.Range("P46").Value = 'Kan ikke legge til %1$d bilder'
$ws.Range("Q46").Value = 'Kan %1$d afbeeldingen niet toevoegen'
$ws.Range("Y46").Value = '%1$d slik ni bilo dodanih.'
$ws.Range("Y47").Value = '%1$d sliki nista bili dodani.'
$ws.Range("C84").Value = 'Delete file \''%1$s\''?'
$ws.Range("G84").Value = 'Datei \''%1$s\'' löschen?'
$ws.Range("J84").Value = 'Supprimer le fichier \''%1$s\''?'
$ws.Range("M84").Value = 'Eliminare il file \’%1$s\’?'
$ws.Range("Q84").Value = 'Bestand \''%1$s\'' verwijderen?'
$ws.Range("C93").Value = 'Deleted %1$d note'
$ws.Range("J93").Value = '%1$d note supprimée'
$ws.Range("M93").Value = 'Eliminata %1$d nota'
$ws.Range("Q93").Value = 'Verwijderde %1$d notitie'
$ws.Range("C94").Value = 'Deleted %1$d notes'
$ws.Range("J94").Value = '%1$d notes supprimées'
$ws.Range("M94").Value = 'Eliminate %1$d note'
$ws.Range("Q94").Value = 'Verwijderde %1$d notities'
$ws.Range("C150").Value = 'Imported %1$s Note'
$ws.Range("J150").Value = '%1$s note importée'
$ws.Range("M150").Value = 'Importata %1$s nota'
$ws.Range("Q150").Value = 'Geïmporteerde %1$s Notitie'
$ws.Range("C151").Value = 'Imported %1$s Notes'
$ws.Range("J151").Value = '%1$s notes importées'
$ws.Range("M151").Value = 'Importate %1$s note'
$ws.Range("Q151").Value = 'Geïmporteerde %1$s Notities'
$ws.Range("C190").Value = '%1$d more'
$ws.Range("G190").Value = '%1$d mehr'
$ws.Range("J190").Value = '%1$d de plus'
$ws.Range("M190").Value = '…ancora %1$d'
$ws.Range("Q190").Value = '%1$d meer'
$ws.Range("C193").Value = '…%1$d more file'
$ws.Range("G193").Value = '…%1$d weitere Datei'
$ws.Range("J193").Value = '…et %1$d fichier '
$ws.Range("M193").Value = '…%1$d altro file'
$ws.Range("Q193").Value = '…%1$d ander bestand'
$ws.Range("C194").Value = '…%1$d more files'
$ws.Range("G194").Value = '…%1$d weitere Dateien'
$ws.Range("J194").Value = '…et %1$d fichiers'
$ws.Range("M194").Value = '…altri %1$d file'
$ws.Range("Q194").Value = '…%1$d andere bestanden'
$ws.Range("C230").Value = 'Restored %1$d note'
$ws.Range("J230").Value = '%1$d note restaurée'
$ws.Range("M230").Value = 'Ripristinata %1$d nota'
$ws.Range("Q230").Value = 'Herstelde %1$d notitie'
$ws.Range("C231").Value = 'Restored %1$d notes'
$ws.Range("J231").Value = '%1$d notes restaurées'
$ws.Range("M231").Value = 'Ripristinate %1$d note'
$ws.Range("Q231").Value = 'Herstelde %1$d notities'
$ws.Range("C268").Value = 'Unarchived %1$d note'
$ws.Range("J268").Value = '%1$d note désarchivée'
$ws.Range("M268").Value = 'Annullata archiviazione di %1$d nota'
$ws.Range("Q268").Value = 'De-gearchiveerde %1$d notitie'
$ws.Range("C269").Value = 'Unarchived %1$d notes'
$ws.Range("J269").Value = '%1$d notes désarchivées'
$ws.Range("M269").Value = 'Annullata archiviazione di %1$d note'
$ws.Range("Q269").Value = 'De-gearchiveerde %1$d notities'
